$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 11-18 first (they are removed entirely in the new data set)
$ws.Range("A11:D18").EntireRow.Delete()

# Ensure columns A, C, D are treated as text so purely-numeric looking
# values (e.g. "130", "6562") are NOT auto-converted to numbers.
$ws.Range("A2:A10").NumberFormat = "@"
$ws.Range("C2:C10").NumberFormat = "@"
$ws.Range("D2:D10").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "130, 423, 780, 1073"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "130"
$ws.Range("D2").Value = "6562"

# Row 3
$ws.Range("A3").Value = "98, 130, 423"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "130"
$ws.Range("D3").Value = "6427"

# Row 4
$ws.Range("A4").Value = "98, 130, 455, 1073"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "130"
$ws.Range("D4").Value = "6489"

# Row 5
$ws.Range("A5").Value = "98, 130, 748, 780, 1073"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "130"
$ws.Range("D5").Value = "6017"

# Row 6
$ws.Range("A6").Value = "130, 423, 748, 1073, SF"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "130"
$ws.Range("D6").Value = "6202"

# Row 7
$ws.Range("A7").Value = "130, 780, 780, 1073"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "780"
$ws.Range("D7").Value = "6359"

# Row 8
$ws.Range("A8").Value = "98, 98, 130, 455, 780"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "130"
$ws.Range("D8").Value = "6727"

# Row 9
$ws.Range("A9").Value = "98, 130, 423, 1073"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "130"
$ws.Range("D9").Value = "6670"

# Row 10
$ws.Range("A10").Value = "98, 98, 455, 780"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "780"
$ws.Range("D10").Value = "6748"
